$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")
$wmitr = $wb.Worksheets.Item("WMITR")

# --- About sheet updates ---
# Source changed from Census/NerdWallet to Congressional Budget Office, with
# additional citation detail rows added beneath it.
$about.Range("B3").Value = "Congressional Budget Office"
$about.Range("B4").Value = 2012
$about.Range("B5").Value = "Effective Marginal Tax Rates for Low- and Moderate-Income Workers"

$about.Hyperlinks.Add($about.Range("B6"), "http://www.cbo.gov/publication/43709")
$about.Range("B6").Style = "Hyperlink"

$about.Range("B7").Value = "Summary, headline of third paragraph"

# Notes section updated to describe the new CBO-based figure.
$about.Range("A10").Value = "This is the average marginal income tax rate faced by workers who earn less than"
$about.Range("A11").Value = "450% of the federal poverty level."

# --- WMITR sheet updates ---
$wmitr.Range("B2").Value = 0.3

Write-Host "edit complete"
